$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (product / pid_pass)
$ws.Range("C2").Value = 0.015
$ws.Range("D2").Value = 0.037
$ws.Range("E2").Value = 0.053
$ws.Range("F2").Value = 0.08400000000000001
$ws.Range("G2").Value = 0.114
$ws.Range("H2").Value = 0.155

# Row 3 (product / pid_recall)
$ws.Range("C3").Value = 0.002285714285714286
$ws.Range("D3").Value = 0.006253174603174603
$ws.Range("E3").Value = 0.01024761904761904
$ws.Range("F3").Value = 0.01538412698412697
$ws.Range("G3").Value = 0.02327420634920634
$ws.Range("H3").Value = 0.03309246031746031

# Row 4 (product / pid_ndcg)
$ws.Range("C4").Value = 0.015
$ws.Range("D4").Value = 0.01353072127397724
$ws.Range("E4").Value = 0.01317972828356509
$ws.Range("F4").Value = 0.01374058675776681
$ws.Range("G4").Value = 0.01700362346505621
$ws.Range("H4").Value = 0.02068832803497481

# Row 5 (video / pid_pass)
$ws.Range("D5").Value = 0.052
$ws.Range("E5").Value = 0.068

# Row 6 (video / pid_recall)
$ws.Range("D6").Value = 0.006795238095238092
$ws.Range("E6").Value = 0.009174999999999994
$ws.Range("F6").Value = 0.01514920634920633
$ws.Range("G6").Value = 0.01754761904761903
$ws.Range("H6").Value = 0.02189325396825395

# Row 7 (video / pid_ndcg)
$ws.Range("D7").Value = 0.01880454152732997
$ws.Range("E7").Value = 0.01645908578346691
$ws.Range("F7").Value = 0.01569074768415001
$ws.Range("G7").Value = 0.01701082202936499
$ws.Range("H7").Value = 0.0190636800964749
